$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 5797.5
$ws.Range("J43").Value = 7063.3335
$ws.Range("L43").Value = 7063.3335
$ws.Range("N43").Value = -7201.3335

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 19010.691
$ws.Range("I32").Value = 10496.619
$ws.Range("J32").Value = 24778.291
$ws.Range("K32").Value = 10496.619
$ws.Range("L32").Value = 24778.291
$ws.Range("M32").Value = -10209.619
$ws.Range("N32").Value = -25352.291

# Row 45
$ws.Range("H45").Value = 2199.2
$ws.Range("I45").Value = 2199.2
$ws.Range("K45").Value = 2199.2
$ws.Range("M45").Value = -1822.2

# Row 88
$ws.Range("H88").Value = 2258.2
$ws.Range("I88").Value = 1947.8334
$ws.Range("J88").Value = 2723.75
$ws.Range("K88").Value = 1947.8334
$ws.Range("L88").Value = 2723.75
$ws.Range("M88").Value = -1541.8334
$ws.Range("N88").Value = -3535.75

# Row 91
$ws.Range("H91").Value = 2258.2
$ws.Range("I91").Value = 1947.8334
$ws.Range("J91").Value = 2723.75
$ws.Range("K91").Value = 1947.8334
$ws.Range("L91").Value = 2723.75
$ws.Range("M91").Value = -543.8334
$ws.Range("N91").Value = -5531.75

# Row 110
$ws.Range("H110").Value = 11099.25
$ws.Range("I110").Value = 11099.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 11099.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -9054.25
$ws.Range("N110").ClearContents()

# Row 132
$ws.Range("H132").Value = 1765.1296
$ws.Range("I132").Value = 1712.52
$ws.Range("K132").Value = 5137.559999999999
$ws.Range("M132").Value = -2607.559999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 8500
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -1753
$ws.Range("N20").Value = -15494

# Row 94
$ws.Range("H94").Value = 504.22223
$ws.Range("I94").Value = 574.6667
$ws.Range("J94").Value = 363.33334
$ws.Range("K94").Value = 574.6667
$ws.Range("L94").Value = 363.33334
$ws.Range("M94").Value = -123.6667
$ws.Range("N94").Value = -1265.33334

# Row 99
$ws.Range("H99").Value = 2535.7144
$ws.Range("I99").Value = 2651
$ws.Range("K99").Value = 2651
$ws.Range("M99").Value = -1153

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 1213.9
$ws.Range("I19").Value = 1452.375
$ws.Range("J19").Value = 260
$ws.Range("K19").Value = 1452.375
$ws.Range("L19").Value = 260
$ws.Range("M19").Value = -1282.375
$ws.Range("N19").Value = -600

# Row 24
$ws.Range("H24").Value = 1213.9
$ws.Range("I24").Value = 1452.375
$ws.Range("J24").Value = 260
$ws.Range("K24").Value = 1452.375
$ws.Range("L24").Value = 260
$ws.Range("M24").Value = -1282.375
$ws.Range("N24").Value = -600

# Row 31
$ws.Range("H31").Value = 4533
$ws.Range("I31").Value = 2464.1428
$ws.Range("J31").Value = 5739.8335
$ws.Range("K31").Value = 2464.1428
$ws.Range("L31").Value = 5739.8335
$ws.Range("M31").Value = -2169.1428
$ws.Range("N31").Value = -6329.8335

# Row 32
$ws.Range("H32").Value = 502000
$ws.Range("I32").Value = 1000000
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 1000000
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -999684
$ws.Range("N32").Value = -4632

# Row 34
$ws.Range("H34").Value = 4533
$ws.Range("I34").Value = 2464.1428
$ws.Range("J34").Value = 5739.8335
$ws.Range("K34").Value = 2464.1428
$ws.Range("L34").Value = 5739.8335
$ws.Range("M34").Value = -2262.1428
$ws.Range("N34").Value = -6143.8335

# Row 62
$ws.Range("H62").Value = 32463.143
$ws.Range("I62").Value = 4021.4443
$ws.Range("K62").Value = 4021.4443
$ws.Range("M62").Value = -3397.4443

# Row 65
$ws.Range("H65").Value = 32463.143
$ws.Range("I65").Value = 4021.4443
$ws.Range("K65").Value = 20107.2215
$ws.Range("M65").Value = -16987.2215

# Row 109
$ws.Range("H109").Value = 13058.667
$ws.Range("J109").Value = 13058.667
$ws.Range("L109").Value = 13058.667
$ws.Range("N109").Value = -15138.667

# Row 122
$ws.Range("H122").Value = 8158.8335
$ws.Range("J122").Value = 7485.6665
$ws.Range("L122").Value = 22456.9995
$ws.Range("N122").Value = -27356.9995

# Row 132
$ws.Range("H132").Value = 2955.2593
$ws.Range("I132").Value = 2696.7
$ws.Range("J132").Value = 3694
$ws.Range("K132").Value = 8090.099999999999
$ws.Range("L132").Value = 11082
$ws.Range("M132").Value = -5560.099999999999
$ws.Range("N132").Value = -16142

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2201489.5
$ws.Range("J4").Value = 3249
$ws.Range("L4").Value = 9747
$ws.Range("N4").Value = -9971

# Row 12
$ws.Range("H12").Value = 32.384617
$ws.Range("I12").Value = 19.142857
$ws.Range("K12").Value = 57.428571
$ws.Range("M12").Value = 115.571429

# Row 75
$ws.Range("H75").Value = 730.375
$ws.Range("J75").Value = 777.2
$ws.Range("L75").Value = 2331.6
$ws.Range("N75").Value = -4327.6

# Row 78
$ws.Range("H78").Value = 730.375
$ws.Range("J78").Value = 777.2
$ws.Range("L78").Value = 6994.8
$ws.Range("N78").Value = -16978.8

# Row 117
$ws.Range("H117").Value = 1859.5834
$ws.Range("J117").Value = 2212.7778
$ws.Range("L117").Value = 6638.3334
$ws.Range("N117").Value = -13522.3334

# Row 119
$ws.Range("H119").Value = 3999
$ws.Range("I119").Value = 3999
$ws.Range("K119").Value = 11997
$ws.Range("M119").Value = -7159

# Row 120
$ws.Range("H120").Value = 14370.952
$ws.Range("I120").Value = 4596.6665
$ws.Range("K120").Value = 13789.9995
$ws.Range("M120").Value = -8951.999500000002

# Row 129
$ws.Range("H129").Value = 3930.111
$ws.Range("I129").Value = 4949.25
$ws.Range("J129").Value = 3114.8
$ws.Range("K129").Value = 14847.75
$ws.Range("L129").Value = 9344.400000000001
$ws.Range("M129").Value = -9847.75
$ws.Range("N129").Value = -19344.4

$ws = $wb.Worksheets.Item("GSM")
# Row 135
$ws.Range("H135").Value = 5049999.5
$ws.Range("J135").Value = 5049999.5
$ws.Range("L135").Value = 5049999.5
$ws.Range("N135").Value = -5060139.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 6307.9443
$ws.Range("I16").Value = 4436.3335
$ws.Range("K16").Value = 4436.3335
$ws.Range("M16").Value = -4266.3335

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# Row 105
$ws.Range("H105").Value = 36969
$ws.Range("J105").Value = 36969
$ws.Range("L105").Value = 36969
$ws.Range("N105").Value = -43957

# Row 122
$ws.Range("H122").Value = 3699
$ws.Range("I122").Value = 3699
$ws.Range("J122").Value = 3699
$ws.Range("K122").Value = 11097
$ws.Range("L122").Value = 11097
$ws.Range("M122").Value = -8647
$ws.Range("N122").Value = -15997

# Row 132
$ws.Range("H132").Value = 3449.5151
$ws.Range("I132").Value = 2297.1304
$ws.Range("K132").Value = 6891.3912
$ws.Range("M132").Value = -4361.3912

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 7000
$ws.Range("J5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("N5").Value = -7224

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 45
$ws.Range("H45").Value = 19018
$ws.Range("I45").Value = 21963.334
$ws.Range("J45").Value = 14600
$ws.Range("K45").Value = 21963.334
$ws.Range("L45").Value = 14600
$ws.Range("M45").Value = -21472.334
$ws.Range("N45").Value = -15582

# Row 113
$ws.Range("H113").Value = 3744.9092
$ws.Range("I113").Value = 2850.25
$ws.Range("J113").Value = 4256.143
$ws.Range("K113").Value = 8550.75
$ws.Range("L113").Value = 12768.429
$ws.Range("M113").Value = -6380.75
$ws.Range("N113").Value = -17108.429

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 126
$ws.Range("H126").Value = 79739
$ws.Range("I126").Value = 112899.89
$ws.Range("K126").Value = 338699.67
$ws.Range("M126").Value = -336229.67
